# Add a new worksheet "TC004" (a new test case), positioned right before
# "Sheet3", populated like TC002/TC003's header + first data row, plus two
# styled-but-empty cells (D4, D5) mirroring TC003's "extra colour" column.

$wb = $excel.ActiveWorkbook

$wsTC002  = $wb.Worksheets.Item("TC002")
$wsTC003  = $wb.Worksheets.Item("TC003")
$wsSheet3 = $wb.Worksheets.Item("Sheet3")

# Insert the new sheet immediately before "Sheet3" so the tab order becomes
# TC002, TC003, TC004, Sheet3.
$wsTC004 = $wb.Worksheets.Add($wsSheet3)
$wsTC004.Name = "TC004"

# Copy header row + first data row (values + styles) from TC002.
$wsTC002.Range("A1:D2").Copy($wsTC004.Range("A1"))

# Reproduce the styled-but-empty D4 / D5 cells (style used by TC003's last
# "color" column) without carrying over any value.
$wsTC003.Range("D4").Copy($wsTC004.Range("D4"))
$wsTC003.Range("D4").Copy($wsTC004.Range("D5"))
$wsTC004.Range("D4:D5").ClearContents()

# Match the author's selection on the new sheet.
$wsTC004.Range("B2").Select()
